# Weekly data refresh: insert the newest price record as the new row 478
# (pushing the existing rows 478-580 down to 479-581) on the single data
# sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 478..580 down one position, creating an empty row 478.
$ws.Rows.Item(478).Insert()

# Populate the newly inserted row 478 with the new daily record.
$ws.Cells.Item(478, 1).Value = 8
$ws.Cells.Item(478, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(478, 3).Value = "Coquimbo"
$ws.Cells.Item(478, 4).Value = 45173
$ws.Cells.Item(478, 5).Value = 4
$ws.Cells.Item(478, 6).Value = 100112032
$ws.Cells.Item(478, 7).Value = "Zapallo italiano"
$ws.Cells.Item(478, 8).Value = "Sin especificar"
$ws.Cells.Item(478, 9).Value = "Primera"
$ws.Cells.Item(478, 10).Value = 300
$ws.Cells.Item(478, 11).Value = 12000
$ws.Cells.Item(478, 12).Value = 13000
$ws.Cells.Item(478, 13).Value = 12500
$ws.Cells.Item(478, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(478, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(478, 16).Value = 250
$ws.Cells.Item(478, 17).Value = 50
$ws.Cells.Item(478, 18).Value = "Hortaliza"
